$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.1017261348290648
$ws.Cells.Item(2, 8).Value = 21.29352360608017
$ws.Cells.Item(2, 9).Value = -11.93301775704541
$ws.Cells.Item(3, 7).Value = 0.1040845169748119
$ws.Cells.Item(3, 8).Value = -10.54628287640636
$ws.Cells.Item(4, 7).Value = -0.02524906764944638
$ws.Cells.Item(4, 8).Value = -63.68629722668907
$ws.Cells.Item(5, 7).Value = -0.006848189521895681
$ws.Cells.Item(5, 8).Value = 90.46899099537985
$ws.Cells.Item(6, 7).Value = 0.04232963056391978
$ws.Cells.Item(6, 8).Value = 19.08521192960232
$ws.Cells.Item(7, 7).Value = 0.05205215838819094
$ws.Cells.Item(7, 8).Value = 156.4765026849545
$ws.Cells.Item(8, 7).Value = -0.1410076249325467
$ws.Cells.Item(8, 8).Value = 0.06397814148557507
$ws.Cells.Item(9, 7).Value = -0.0931384134039961
$ws.Cells.Item(9, 8).Value = 31.95964083404833
$ws.Cells.Item(10, 7).Value = -0.08140593978849425
$ws.Cells.Item(10, 8).Value = 24.64305498083406
$ws.Cells.Item(11, 7).Value = -0.06838143202799279
$ws.Cells.Item(11, 8).Value = -2.766075379183044
$ws.Cells.Item(12, 7).Value = -0.3062872950165193
$ws.Cells.Item(12, 8).Value = 26.23073044155044
$ws.Cells.Item(13, 7).Value = -0.4799252170516315
$ws.Cells.Item(13, 8).Value = -6.911532927574017
$ws.Cells.Item(14, 7).Value = -0.05834477140278677
$ws.Cells.Item(14, 8).Value = -14.99747906855995
$ws.Cells.Item(15, 7).Value = 0.03166751905520186
$ws.Cells.Item(15, 8).Value = 138.2960074802114
$ws.Cells.Item(16, 7).Value = 0.085056720387856
$ws.Cells.Item(16, 8).Value = -41.75102318890255
$ws.Cells.Item(17, 7).Value = 0.1513522822674456
$ws.Cells.Item(17, 8).Value = 23.4304754316399
$ws.Cells.Item(18, 7).Value = 0.1377588984371945
$ws.Cells.Item(18, 8).Value = -0.5043308002193144
$ws.Cells.Item(19, 7).Value = 0.1245967731651461
$ws.Cells.Item(19, 8).Value = 30.59535089951397
$ws.Cells.Item(20, 7).Value = 0.0238422105451502
$ws.Cells.Item(20, 8).Value = -6.990083304686983
$ws.Cells.Item(21, 7).Value = 0.06612753755988257
$ws.Cells.Item(21, 8).Value = -12.0259707508041
$ws.Cells.Item(24, 7).Value = 0.08597868876261189
$ws.Cells.Item(24, 8).Value = -14.39498310484945
$ws.Cells.Item(25, 7).Value = 0.190992809582956
$ws.Cells.Item(25, 8).Value = 26.02915060973912
$ws.Cells.Item(26, 7).Value = 0.08153210251672185
$ws.Cells.Item(26, 8).Value = 3.042527400864787
$ws.Cells.Item(27, 7).Value = 0.08042933580459205
$ws.Cells.Item(27, 8).Value = -19.45195335169363
$ws.Cells.Item(28, 7).Value = -0.2221245310133198
$ws.Cells.Item(28, 8).Value = -4.225497791813523
$ws.Cells.Item(29, 7).Value = -0.2131768830151092
$ws.Cells.Item(29, 8).Value = -3.83339077250234
$ws.Cells.Item(30, 7).Value = 0.04760835927444473
$ws.Cells.Item(30, 8).Value = 7.878058481954279
$ws.Cells.Item(31, 7).Value = 0.03160195667561526
$ws.Cells.Item(31, 8).Value = 20.00843976363108
$ws.Cells.Item(32, 7).Value = 0.09093899487674141
$ws.Cells.Item(32, 8).Value = -4.234932962798648
$ws.Cells.Item(33, 7).Value = 0.1285008660690005
$ws.Cells.Item(33, 8).Value = 23.60025565536792
$ws.Cells.Item(34, 7).Value = 0.03648537994083615
$ws.Cells.Item(34, 8).Value = -21.41461503317809
$ws.Cells.Item(35, 7).Value = 0.01119523783701483
$ws.Cells.Item(35, 8).Value = 47.76212151055888
$ws.Cells.Item(36, 7).Value = 0.05069934834570633
$ws.Cells.Item(36, 8).Value = -12.1880068058324
$ws.Cells.Item(37, 7).Value = 0.08710397700073065
$ws.Cells.Item(37, 8).Value = 23.85715042662389
$ws.Cells.Item(38, 7).Value = 0.009988417534054259
$ws.Cells.Item(38, 8).Value = -80.93318997592159
$ws.Cells.Item(39, 7).Value = 0.02641625155537421
$ws.Cells.Item(39, 8).Value = 27.38502214990844
$ws.Cells.Item(40, 7).Value = 0.004009011956253929
$ws.Cells.Item(40, 8).Value = 147.2148969967322
$ws.Cells.Item(41, 7).Value = 0.04438057546661396
$ws.Cells.Item(41, 8).Value = 25.52539517991285
$ws.Cells.Item(42, 7).Value = 0.1455590530514656
$ws.Cells.Item(42, 8).Value = 8.881196307085215
$ws.Cells.Item(43, 7).Value = 0.1504298857027398
$ws.Cells.Item(43, 8).Value = 0.9756926944486395
$ws.Cells.Item(44, 7).Value = -0.01154866942518077
$ws.Cells.Item(44, 8).Value = -35.68342638788887
$ws.Cells.Item(45, 7).Value = 0.009738881891338013
$ws.Cells.Item(45, 8).Value = 188.6976464493591
$ws.Cells.Item(46, 7).Value = 0.01332111482796823
$ws.Cells.Item(46, 8).Value = 504.5225969286569
$ws.Cells.Item(47, 7).Value = -0.03985422012966428
$ws.Cells.Item(47, 8).Value = -329.5338829767937
$ws.Cells.Item(48, 7).Value = 0.0732223042245497
$ws.Cells.Item(48, 8).Value = 45.65193562271457
$ws.Cells.Item(49, 7).Value = 0.06192338551369429
$ws.Cells.Item(49, 8).Value = -6.268692474162636
$ws.Cells.Item(50, 7).Value = 0.1523186800358433
$ws.Cells.Item(50, 8).Value = -5.543302039630862
$ws.Cells.Item(51, 7).Value = 0.1662078747793498
$ws.Cells.Item(51, 8).Value = -2.867712272219629
$ws.Cells.Item(52, 7).Value = -0.1695033119434113
$ws.Cells.Item(52, 8).Value = -5.654960845842983
$ws.Cells.Item(53, 7).Value = -0.1590418214801917
$ws.Cells.Item(53, 8).Value = -26.16909434673128
$ws.Cells.Item(54, 7).Value = 0.1195044352546354
$ws.Cells.Item(54, 8).Value = 27.5089427671595
$ws.Cells.Item(55, 7).Value = 0.1078362300526788
$ws.Cells.Item(55, 8).Value = -4.63785104557739
$ws.Cells.Item(56, 7).Value = -0.01650360495763801
$ws.Cells.Item(56, 8).Value = -126.0528897223899
$ws.Cells.Item(57, 7).Value = -0.01693888215723816
$ws.Cells.Item(57, 8).Value = 25.92560306256595
$ws.Cells.Item(58, 7).Value = 0.05696028736865899
$ws.Cells.Item(58, 8).Value = 1.014959006096574
$ws.Cells.Item(59, 7).Value = 0.07296475805731102
$ws.Cells.Item(59, 8).Value = 1.594031282660871
$ws.Cells.Item(60, 7).Value = 0.04719781368694941
$ws.Cells.Item(60, 8).Value = -32.54945395581765
$ws.Cells.Item(61, 7).Value = 0.07077194594613298
$ws.Cells.Item(61, 8).Value = 48.9095994113699
$ws.Cells.Item(62, 7).Value = 0.07611680605606318
$ws.Cells.Item(62, 8).Value = 4.326294373761176
$ws.Cells.Item(63, 7).Value = 0.0756659439096507
$ws.Cells.Item(63, 8).Value = 15.70865400003319
$ws.Cells.Item(64, 7).Value = -0.04186531977830304
$ws.Cells.Item(64, 8).Value = -1.086096123162964
$ws.Cells.Item(65, 7).Value = 0.0272336423276875
$ws.Cells.Item(65, 8).Value = 155.2141021081198
$ws.Cells.Item(66, 7).Value = 0.02878263240873565
$ws.Cells.Item(66, 8).Value = 52.01579060647462
$ws.Cells.Item(67, 7).Value = 0.03113408603254171
$ws.Cells.Item(67, 8).Value = 19.06012211986882
$ws.Cells.Item(68, 7).Value = -0.01431713002426085
$ws.Cells.Item(68, 8).Value = -2611.982667491071
$ws.Cells.Item(69, 7).Value = 0.01496830221565967
$ws.Cells.Item(69, 8).Value = 215.9227584676107
$ws.Cells.Item(70, 7).Value = -0.03355111833996965
$ws.Cells.Item(70, 8).Value = -22.23489900008101
$ws.Cells.Item(71, 7).Value = -0.03885132597881324
$ws.Cells.Item(71, 8).Value = 29.47912484346245
$ws.Cells.Item(72, 7).Value = -0.1265992073410618
$ws.Cells.Item(72, 8).Value = 14.65268746056415
$ws.Cells.Item(73, 7).Value = -0.1569735785354142
$ws.Cells.Item(73, 8).Value = -8.414777493084422
$ws.Cells.Item(74, 7).Value = 0.1218881767594388
$ws.Cells.Item(74, 8).Value = -3.302572282404267
$ws.Cells.Item(75, 7).Value = 0.160546477925136
$ws.Cells.Item(75, 8).Value = 18.77201268500899
$ws.Cells.Item(76, 7).Value = -0.06542549730318846
$ws.Cells.Item(76, 8).Value = -89.98084010277421
$ws.Cells.Item(77, 7).Value = -0.03175318904640758
$ws.Cells.Item(77, 8).Value = 31.25082023802416
$ws.Cells.Item(78, 7).Value = 0.09199803682596355
$ws.Cells.Item(78, 8).Value = -0.1882374474335531
$ws.Cells.Item(79, 7).Value = 0.09498114591721894
$ws.Cells.Item(79, 8).Value = -1.578459728478865
$ws.Cells.Item(80, 7).Value = -0.1886796654908553
$ws.Cells.Item(80, 8).Value = -16.17286177856309
$ws.Cells.Item(81, 7).Value = -0.1891600139371402
$ws.Cells.Item(81, 8).Value = 12.60736358712254
$ws.Cells.Item(82, 7).Value = 0.1623815879743861
$ws.Cells.Item(82, 8).Value = 17.03491368822547
$ws.Cells.Item(83, 7).Value = 0.203840514984346
$ws.Cells.Item(83, 8).Value = 23.82518315064771
$ws.Cells.Item(84, 7).Value = 0.06732362981211748
$ws.Cells.Item(84, 8).Value = 381.0539280694853
$ws.Cells.Item(85, 7).Value = 0.002897458319204544
$ws.Cells.Item(85, 8).Value = -87.20157777311987
